# Generate Report for Handoff
# Adds two new files (32177ba9-... and 8041dc99-...) to the localization
# status report, on the Overview sheet and on each language sheet
# (zh-cn, de-de). The previously-last row (b9e30e55-...) is pushed down
# to make room, keeping its original values.
#
# NB: deleting a hyperlink via any single cell's Hyperlinks collection
# clears *all* hyperlinks on the sheet in this runtime, so every sheet's
# hyperlinks are cleared up-front and then fully re-created (old + new)
# in final left-to-right / top-to-bottom order, matching how Excel
# numbers relationship ids when a user edits such a sheet interactively.

$wb = $excel.ActiveWorkbook

$repoBase        = "https://github.com/OpenLocalizationTest/oltest/blob"
$handoffOrgBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Hyperlinks.Delete()

$ws1.Cells.Item(3,2).Value = "Ready for handoff"
$ws1.Cells.Item(3,3).Value = "Ready for handoff"
$ws1.Cells.Item(3,4).Value = "2016-25-14 02:25:30"

$ws1.Cells.Item(4,2).Value = "Ready for handoff"
$ws1.Cells.Item(4,3).Value = "Ready for handoff"
$ws1.Cells.Item(4,4).Value = "2016-25-14 02:25:30"

$ws1.Cells.Item(5,2).Value = "Ready for handoff"
$ws1.Cells.Item(5,3).Value = "Ready for handoff"
$ws1.Cells.Item(5,4).Value = "2016-24-14 02:24:14"

$ws1.Hyperlinks.Add($ws1.Cells.Item(2,1), "$repoBase/ce321fe7253297d258f2156f2d02c6b86a5100af/e2e/d3b6309e-e5e3-4b7f-a0a3-733545ef4820.md", "", "", "d3b6309e-e5e3-4b7f-a0a3-733545ef4820.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(3,1), "$repoBase/29bb3d245767ab20c8d2bafdb653d824ddcde021/e2e/32177ba9-78a3-42cf-a90e-85b40a7a2e73.md", "", "", "32177ba9-78a3-42cf-a90e-85b40a7a2e73.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(4,1), "$repoBase/3b694909c966ce442d3b2aa1d12523216c5ab3ac/e2e/8041dc99-f239-4c18-830e-179da63b0719.md", "", "", "8041dc99-f239-4c18-830e-179da63b0719.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(5,1), "$repoBase/f373e3f880e833c536ef4a092a050b3ef0d39282/e2e/b9e30e55-7b5a-4e23-aaf2-5a8576e57075.md", "", "", "b9e30e55-7b5a-4e23-aaf2-5a8576e57075.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Delete()

$ws2.Cells.Item(3,3).Value = "Ready for handoff"
$ws2.Cells.Item(3,5).Value = "2016-03-14 02:25:28"
$ws2.Cells.Item(3,8).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(3,9).Value = "Include"

$ws2.Cells.Item(4,3).Value = "Ready for handoff"
$ws2.Cells.Item(4,5).Value = "2016-03-14 02:25:28"
$ws2.Cells.Item(4,8).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(4,9).Value = "Include"

$ws2.Cells.Item(5,3).Value = "Ready for handoff"
$ws2.Cells.Item(5,5).Value = "2016-03-14 02:24:12"
$ws2.Cells.Item(5,8).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(5,9).Value = "Include"

$ws2.Hyperlinks.Add($ws2.Cells.Item(2,1), "$repoBase/ce321fe7253297d258f2156f2d02c6b86a5100af/e2e/d3b6309e-e5e3-4b7f-a0a3-733545ef4820.md", "", "", "d3b6309e-e5e3-4b7f-a0a3-733545ef4820.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(2,2), "$repoBase/ce321fe7253297d258f2156f2d02c6b86a5100af/e2e/d3b6309e-e5e3-4b7f-a0a3-733545ef4820.md", "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(2,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/52f48e73835f8c1aa021049292d792f40f805301/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d3b6309e-e5e3-4b7f-a0a3-733545ef4820.d86ef5636691266e9b48cb3d305636d389f867b2.zh-cn.xlf", "", "", "d3b6309e-e5e3-4b7f-a0a3-733545ef4820.d86ef5636691266e9b48cb3d305636d389f867b2.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(2,6), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/dc2aaac436c17f46401ce346f1eca87aee9508f9/e2e/d3b6309e-e5e3-4b7f-a0a3-733545ef4820.md", "", "", "d3b6309e-e5e3-4b7f-a0a3-733545ef4820.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(2,7), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/02e00c26c1ee40bf475da8f4201c5b078977488e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d3b6309e-e5e3-4b7f-a0a3-733545ef4820.d86ef5636691266e9b48cb3d305636d389f867b2.zh-cn.xlf", "", "", "d3b6309e-e5e3-4b7f-a0a3-733545ef4820.d86ef5636691266e9b48cb3d305636d389f867b2.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Cells.Item(3,1), "$repoBase/29bb3d245767ab20c8d2bafdb653d824ddcde021/e2e/32177ba9-78a3-42cf-a90e-85b40a7a2e73.md", "", "", "32177ba9-78a3-42cf-a90e-85b40a7a2e73.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(3,2), "$repoBase/29bb3d245767ab20c8d2bafdb653d824ddcde021/e2e/32177ba9-78a3-42cf-a90e-85b40a7a2e73.md", "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(3,4), "$handoffOrgBase/41e0d616becbbecd16aabcf92b762057331b18b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/32177ba9-78a3-42cf-a90e-85b40a7a2e73.3f174ce30e4fc1518b21b3fee4290539e59c96d4.zh-cn.xlf", "", "", "32177ba9-78a3-42cf-a90e-85b40a7a2e73.3f174ce30e4fc1518b21b3fee4290539e59c96d4.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Cells.Item(4,1), "$repoBase/3b694909c966ce442d3b2aa1d12523216c5ab3ac/e2e/8041dc99-f239-4c18-830e-179da63b0719.md", "", "", "8041dc99-f239-4c18-830e-179da63b0719.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(4,2), "$repoBase/3b694909c966ce442d3b2aa1d12523216c5ab3ac/e2e/8041dc99-f239-4c18-830e-179da63b0719.md", "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(4,4), "$handoffOrgBase/afd5036e68b0c0cfa778b6aad0c81dfc1b78900c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8041dc99-f239-4c18-830e-179da63b0719.686baace7255c656eb06a8fefc835a09dd9116e4.zh-cn.xlf", "", "", "8041dc99-f239-4c18-830e-179da63b0719.686baace7255c656eb06a8fefc835a09dd9116e4.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Cells.Item(5,1), "$repoBase/f373e3f880e833c536ef4a092a050b3ef0d39282/e2e/b9e30e55-7b5a-4e23-aaf2-5a8576e57075.md", "", "", "b9e30e55-7b5a-4e23-aaf2-5a8576e57075.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(5,2), "$repoBase/f373e3f880e833c536ef4a092a050b3ef0d39282/e2e/b9e30e55-7b5a-4e23-aaf2-5a8576e57075.md", "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(5,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/072accb88c9e9650b0f05671655245dfc31adbf5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b9e30e55-7b5a-4e23-aaf2-5a8576e57075.c6b89e88f0c8f40cb01622f1eee960e6aa8e0816.zh-cn.xlf", "", "", "b9e30e55-7b5a-4e23-aaf2-5a8576e57075.c6b89e88f0c8f40cb01622f1eee960e6aa8e0816.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Hyperlinks.Delete()

$ws3.Cells.Item(3,3).Value = "Ready for handoff"
$ws3.Cells.Item(3,5).Value = "2016-03-14 02:25:30"
$ws3.Cells.Item(3,8).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(3,9).Value = "Include"

$ws3.Cells.Item(4,3).Value = "Ready for handoff"
$ws3.Cells.Item(4,5).Value = "2016-03-14 02:25:30"
$ws3.Cells.Item(4,8).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(4,9).Value = "Include"

$ws3.Cells.Item(5,3).Value = "Ready for handoff"
$ws3.Cells.Item(5,5).Value = "2016-03-14 02:24:14"
$ws3.Cells.Item(5,8).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(5,9).Value = "Include"

$ws3.Hyperlinks.Add($ws3.Cells.Item(2,1), "$repoBase/ce321fe7253297d258f2156f2d02c6b86a5100af/e2e/d3b6309e-e5e3-4b7f-a0a3-733545ef4820.md", "", "", "d3b6309e-e5e3-4b7f-a0a3-733545ef4820.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(2,2), "$repoBase/ce321fe7253297d258f2156f2d02c6b86a5100af/e2e/d3b6309e-e5e3-4b7f-a0a3-733545ef4820.md", "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(2,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/175037f07b522ec0549f16047a3adbe7895fe736/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d3b6309e-e5e3-4b7f-a0a3-733545ef4820.d86ef5636691266e9b48cb3d305636d389f867b2.de-de.xlf", "", "", "d3b6309e-e5e3-4b7f-a0a3-733545ef4820.d86ef5636691266e9b48cb3d305636d389f867b2.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(2,6), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/d48bf3d2b8b3338084833634364afa78f1e51d0a/e2e/d3b6309e-e5e3-4b7f-a0a3-733545ef4820.md", "", "", "d3b6309e-e5e3-4b7f-a0a3-733545ef4820.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(2,7), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ec0bdbbb6751b771d0120e28a3ad952a7ff1899a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d3b6309e-e5e3-4b7f-a0a3-733545ef4820.d86ef5636691266e9b48cb3d305636d389f867b2.de-de.xlf", "", "", "d3b6309e-e5e3-4b7f-a0a3-733545ef4820.d86ef5636691266e9b48cb3d305636d389f867b2.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Cells.Item(3,1), "$repoBase/29bb3d245767ab20c8d2bafdb653d824ddcde021/e2e/32177ba9-78a3-42cf-a90e-85b40a7a2e73.md", "", "", "32177ba9-78a3-42cf-a90e-85b40a7a2e73.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(3,2), "$repoBase/29bb3d245767ab20c8d2bafdb653d824ddcde021/e2e/32177ba9-78a3-42cf-a90e-85b40a7a2e73.md", "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(3,4), "$handoffOrgBase/0d4ec8653f65911534189e1aefd07b2147dfec8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/32177ba9-78a3-42cf-a90e-85b40a7a2e73.3f174ce30e4fc1518b21b3fee4290539e59c96d4.de-de.xlf", "", "", "32177ba9-78a3-42cf-a90e-85b40a7a2e73.3f174ce30e4fc1518b21b3fee4290539e59c96d4.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Cells.Item(4,1), "$repoBase/3b694909c966ce442d3b2aa1d12523216c5ab3ac/e2e/8041dc99-f239-4c18-830e-179da63b0719.md", "", "", "8041dc99-f239-4c18-830e-179da63b0719.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(4,2), "$repoBase/3b694909c966ce442d3b2aa1d12523216c5ab3ac/e2e/8041dc99-f239-4c18-830e-179da63b0719.md", "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(4,4), "$handoffOrgBase/87147391a5f461eab843b20758b57e9c070574b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8041dc99-f239-4c18-830e-179da63b0719.686baace7255c656eb06a8fefc835a09dd9116e4.de-de.xlf", "", "", "8041dc99-f239-4c18-830e-179da63b0719.686baace7255c656eb06a8fefc835a09dd9116e4.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Cells.Item(5,1), "$repoBase/f373e3f880e833c536ef4a092a050b3ef0d39282/e2e/b9e30e55-7b5a-4e23-aaf2-5a8576e57075.md", "", "", "b9e30e55-7b5a-4e23-aaf2-5a8576e57075.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(5,2), "$repoBase/f373e3f880e833c536ef4a092a050b3ef0d39282/e2e/b9e30e55-7b5a-4e23-aaf2-5a8576e57075.md", "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(5,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f40e985edd161dd8a98d3993d765f86eb08305a7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b9e30e55-7b5a-4e23-aaf2-5a8576e57075.c6b89e88f0c8f40cb01622f1eee960e6aa8e0816.de-de.xlf", "", "", "b9e30e55-7b5a-4e23-aaf2-5a8576e57075.c6b89e88f0c8f40cb01622f1eee960e6aa8e0816.de-de.xlf") | Out-Null
